# Auto-generated Word COM-interop script applying the moveset-document edit.
$d = $word.ActiveDocument

# -------------------------------------------------------------------------
# 1) Plain-text corrections via Find/Replace (wdReplaceAll literal matches).
# -------------------------------------------------------------------------
function Replace-ExactText($oldText, $newText) {
    $range = $word.ActiveDocument.Content
    $ok = $range.Find.Execute($oldText, $true, $true, $false, $false, $false, $true, 1, $false, $newText, 2)
    if (-not $ok) {
        throw "Find/Replace failed for: $oldText"
    }
}

Replace-ExactText 'j.B- Far reaching, horizontal kick.' 'j.B- Jumping kick at a low angle.'
Replace-ExactText 'j.C- Jumping kick at a low angle.' 'j.C-  Far reaching, horizontal kick.'
Replace-ExactText '214(E) -- Vine Transit- Envelop yourself in vines and travel underground (Zato BTL)' '214(E) -- Vine Transit- Envelop yourself in vines and travel underground (Similar to Zato’s Break The Law from Guilty Gear)'
Replace-ExactText 'Element Super(236A/B/C + E) -- Soul Transfusion- Healing field that restores owner HP and decreases enemy HP while they are both in the field. Last 7 seconds.' 'Element Super(236A/B/C + E) -- Soul Transfusion- Healing field that restores owner HP and decreases enemy HP while they are both in the field. Last 6 seconds.'
Replace-ExactText 'Element Super(236A/B/C + E) -- Impending Death- Summon a slow moving, multi-hitting projectile(Similar to Dormammu’s Stalking Flare)' 'Element Super(236A/B/C + E) -- Impending Death- Summon a slow moving, multi-hitting projectile(Similar to Dormammu’s Stalking Flare from Ultimate Marvel vs Capcom 3)'
Replace-ExactText 'Passive: Extra Air Option' 'Passive: Extra Air Movement Option'
Replace-ExactText 'Element Super(236A/B/C + E) -- Time Skip- Quickly rushes forward and attacks the opponent. Projectile invincible. (Johnny’s SDTH)' 'Element Super(236A/B/C + E) -- Time Skip- Quickly rushes forward and attacks the opponent. Projectile invincible. (Similar to Johnny’s Treasure Hunt from Guilty Gear)'

# -------------------------------------------------------------------------
# 2) Picture rename metadata (wp:docPr/@name + pic:cNvPr/@name). These are
#    non-visual drawing properties that Find/Replace cannot reach, and
#    InlineShape has no settable Name property, so each affected picture's
#    enclosing paragraph is rewritten in place via Range.InsertXML with the
#    same XML except for the two name attributes.
# -------------------------------------------------------------------------
function Rename-PictureInParagraph($paraIndex, $oldName, $newName, $newParagraphXml) {
    $para = $word.ActiveDocument.Paragraphs.Item($paraIndex)
    $range = $para.Range
    $range.InsertXML($newParagraphXml)
}

Rename-PictureInParagraph 19 'image19.png' 'image7.png' '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w:rsidR="00000000" w:rsidDel="00000000" w:rsidP="00000000" w:rsidRDefault="00000000" w:rsidRPr="00000000" w14:paraId="00000012"><w:pPr><w:contextualSpacing w:val="0"/><w:rPr/></w:pPr><w:r w:rsidDel="00000000" w:rsidR="00000000" w:rsidRPr="00000000"><w:rPr/><w:drawing><wp:inline distB="114300" distT="114300" distL="114300" distR="114300"><wp:extent cx="1685925" cy="1200150"/><wp:effectExtent b="0" l="0" r="0" t="0"/><wp:docPr id="2" name="image7.png"/><a:graphic><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic><pic:nvPicPr><pic:cNvPr id="0" name="image7.png"/><pic:cNvPicPr preferRelativeResize="0"/></pic:nvPicPr><pic:blipFill><a:blip r:embed="rId14"/><a:srcRect b="16091" l="10747" r="6542" t="11494"/><a:stretch><a:fillRect/></a:stretch></pic:blipFill><pic:spPr><a:xfrm><a:off x="0" y="0"/><a:ext cx="1685925" cy="1200150"/></a:xfrm><a:prstGeom prst="rect"/><a:ln/></pic:spPr></pic:pic></a:graphicData></a:graphic></wp:inline></w:drawing></w:r><w:r w:rsidDel="00000000" w:rsidR="00000000" w:rsidRPr="00000000"><w:rPr><w:rtl w:val="0"/></w:rPr></w:r></w:p>'
Rename-PictureInParagraph 47 'image20.png' 'image12.png' '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w:rsidR="00000000" w:rsidDel="00000000" w:rsidP="00000000" w:rsidRDefault="00000000" w:rsidRPr="00000000" w14:paraId="0000002E"><w:pPr><w:contextualSpacing w:val="0"/><w:rPr/></w:pPr><w:r w:rsidDel="00000000" w:rsidR="00000000" w:rsidRPr="00000000"><w:rPr/><w:drawing><wp:inline distB="114300" distT="114300" distL="114300" distR="114300"><wp:extent cx="1405063" cy="1443038"/><wp:effectExtent b="0" l="0" r="0" t="0"/><wp:docPr id="3" name="image12.png"/><a:graphic><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic><pic:nvPicPr><pic:cNvPr id="0" name="image12.png"/><pic:cNvPicPr preferRelativeResize="0"/></pic:nvPicPr><pic:blipFill><a:blip r:embed="rId18"/><a:srcRect b="0" l="0" r="0" t="0"/><a:stretch><a:fillRect/></a:stretch></pic:blipFill><pic:spPr><a:xfrm><a:off x="0" y="0"/><a:ext cx="1405063" cy="1443038"/></a:xfrm><a:prstGeom prst="rect"/><a:ln/></pic:spPr></pic:pic></a:graphicData></a:graphic></wp:inline></w:drawing></w:r><w:r w:rsidDel="00000000" w:rsidR="00000000" w:rsidRPr="00000000"><w:rPr><w:rtl w:val="0"/></w:rPr></w:r></w:p>'
Rename-PictureInParagraph 87 'image18.png' 'image3.png' '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w:rsidR="00000000" w:rsidDel="00000000" w:rsidP="00000000" w:rsidRDefault="00000000" w:rsidRPr="00000000" w14:paraId="00000056"><w:pPr><w:contextualSpacing w:val="0"/><w:rPr/></w:pPr><w:r w:rsidDel="00000000" w:rsidR="00000000" w:rsidRPr="00000000"><w:rPr/><w:drawing><wp:inline distB="114300" distT="114300" distL="114300" distR="114300"><wp:extent cx="523875" cy="1019175"/><wp:effectExtent b="0" l="0" r="0" t="0"/><wp:docPr id="1" name="image3.png"/><a:graphic><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic><pic:nvPicPr><pic:cNvPr id="0" name="image3.png"/><pic:cNvPicPr preferRelativeResize="0"/></pic:nvPicPr><pic:blipFill><a:blip r:embed="rId32"/><a:srcRect b="0" l="0" r="0" t="0"/><a:stretch><a:fillRect/></a:stretch></pic:blipFill><pic:spPr><a:xfrm><a:off x="0" y="0"/><a:ext cx="523875" cy="1019175"/></a:xfrm><a:prstGeom prst="rect"/><a:ln/></pic:spPr></pic:pic></a:graphicData></a:graphic></wp:inline></w:drawing></w:r><w:r w:rsidDel="00000000" w:rsidR="00000000" w:rsidRPr="00000000"><w:rPr><w:rtl w:val="0"/></w:rPr></w:r></w:p>'

Write-Output 'Edits applied.'
